$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (currently at the very start of
#    the document, around the "ASP.NET Core - API" heading).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Delete the whole paragraph "Partially updating a resource".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Partially updating a resource`r") {
        $p.Range.Delete()
        break
    }
}

# 3. Re-add the "_GoBack" bookmark as a collapsed bookmark right at the end
#    of the "A validation Alternative" paragraph (after its run, before the
#    paragraph mark). A collapsed range placed exactly one position before a
#    paragraph's end confuses Bookmarks.Add in this runtime, so work around
#    it: temporarily append a marker character, drop the (now safely
#    interior) bookmark next to it, then remove the marker again - the
#    bookmark stays put.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "A validation Alternative`r") {
        $endPos = $p.Range.End - 1
        $marker = $d.Range($endPos, $endPos)
        $marker.InsertAfter("X")

        $p2 = $d.Paragraphs.Item($i)
        $bmPos = $p2.Range.End - 2
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        $p3 = $d.Paragraphs.Item($i)
        $xPos = $p3.Range.End - 2
        $xRange = $d.Range($xPos, $xPos + 1)
        $xRange.Delete()
        break
    }
}
